$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Save" column header in H1, copying the header formatting
# from the existing G1 ("sum") header cell so it keeps the bold/centered/
# bordered style (s="1") used by the rest of the header row.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Fill in the "Save" values for each data row (rows 2-10). These are
# literal values from the source data (1 where a save occurred, 0 otherwise).
$saveValues = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 0
    6  = 0
    7  = 0
    8  = 1
    9  = 0
    10 = 1
}

foreach ($row in $saveValues.Keys) {
    $ws.Range("H$row").Value = $saveValues[$row]
}

$excel.CutCopyMode = 0
